$wb = $excel.ActiveWorkbook

$wsShort = $wb.Worksheets.Item("Short Term")
$wsMedium = $wb.Worksheets.Item("Medium Term")

# --- Sheet "Short Term": update existing rows 114-119 ---
$shortUpdates = @{
    114 = @{ B = -2.75; C = -9.76; D = -1.92; E = 7.43; F = 8.66; G = -10.37 }
    115 = @{ B = -2.27; C = -0.94; D = 6.49 }
    116 = @{ B = 6;     D = 3.25 }
    117 = @{ B = 0.97;  D = 6.53 }
    118 = @{ B = 20.07; C = 23.73; D = -11.74 }
    119 = @{ B = -9.5;  C = -6.85; D = -3.42; E = 25.73; F = 21.42; G = -2.37 }
}

foreach ($row in $shortUpdates.Keys) {
    $cols = $shortUpdates[$row]
    foreach ($col in $cols.Keys) {
        $wsShort.Range("$col$row").Value = $cols[$col]
    }
}

# Add new row 120 (copy formatting down from row 119, then set the new values)
$wsShort.Range("A119:G119").Copy($wsShort.Range("A120:G120"))
$wsShort.Range("A120").Value = 45597
$wsShort.Range("B120").Value = 3.07
$wsShort.Range("C120").Value = -4.42
$wsShort.Range("D120").Value = 4.65
$wsShort.Range("E120").Value = 26.54
$wsShort.Range("F120").Value = 27.52
$wsShort.Range("G120").Value = -18.48

# --- Sheet "Medium Term": update existing rows 100-105 ---
$mediumUpdates = @{
    100 = @{ B = 12.1;  C = 7.36;  D = 0.31 }
    101 = @{ B = 8.14;  C = 7.07;  D = -0.8 }
    102 = @{ B = 4.18;  C = 5.94;  D = -0.04 }
    103 = @{ C = 7.94;  D = 2.18 }
    104 = @{ C = 14.51; D = 7.36 }
    105 = @{ B = 23.4;  C = 13.68; D = 9.27 }
}

foreach ($row in $mediumUpdates.Keys) {
    $cols = $mediumUpdates[$row]
    foreach ($col in $cols.Keys) {
        $wsMedium.Range("$col$row").Value = $cols[$col]
    }
}

# Add new row 106 (copy formatting down from row 105, then set the new values)
$wsMedium.Range("A105:D105").Copy($wsMedium.Range("A106:D106"))
$wsMedium.Range("A106").Value = 45597
$wsMedium.Range("B106").Value = 29.94
$wsMedium.Range("C106").Value = 16.76
$wsMedium.Range("D106").Value = 12.09
